$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'57.159.26"
$ws.Range('E2').Value = '  -5.49%  '
$ws.Range('D3').Value = "'2.898.54"
$ws.Range('E3').Value = '  -3.23%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'547.48"
$ws.Range('E5').Value = '  -3.55%  '
$ws.Range('D6').Value = "'124.86"
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = "'0.502"
$ws.Range('E8').Value = '  +0.80%  '
$ws.Range('D9').Value = "'2.885.56"
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').Value = "'0.123"
$ws.Range('E10').Value = '  -7.80%  '
$ws.Range('D11').Value = "'4.66"
$ws.Range('E11').Value = '  -7.84%  '
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').Value = "'0.0000211"
$ws.Range('E13').Value = '  -5.60%  '
$ws.Range('D14').Value = "'31.98"
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = "'3.374.78"
$ws.Range('E16').Value = '  -3.23%  '
$ws.Range('D17').Value = "'2.894.64"
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('D18').Value = "'6.47"
$ws.Range('E18').Value = '  +5.26%  '
$ws.Range('D19').Value = "'57.184.66"
$ws.Range('E19').Value = '  -5.47%  '
$ws.Range('D20').Value = "'402.92"
$ws.Range('E20').Value = '  -6.40%  '
$ws.Range('D21').Value = "'12.75"
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').Value = "'0.667"
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('E23').Value = '  -4.74%  '
$ws.Range('D24').Value = "'12.61"
$ws.Range('E24').Value = '  -2.63%  '
$ws.Range('D25').Value = "'77.41"
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = "'2.45"
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('D29').Value = "'7.18"
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').Value = "'24.55"
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('D32').Value = "'5.92"
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('D33').Value = "'0.0978"
$ws.Range('E33').Value = '  +4.91%  '
$ws.Range('D34').Value = "'0.909"
$ws.Range('E34').Value = '  -3.98%  '
$ws.Range('E35').Value = '  -2.46%  '
$ws.Range('E36').Value = '  -11.71%  '
$ws.Range('D37').Value = "'48.02"
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('D38').Value = "'8.15"
$ws.Range('E38').Value = '  +4.58%  '
$ws.Range('D39').Value = "'0.0₃0633"
$ws.Range('E39').Value = '  -4.25%  '
$ws.Range('D40').Value = "'0.106"
$ws.Range('E40').Value = '  -1.12%  '
$ws.Range('D41').Value = "'0.0336"
$ws.Range('E41').Value = '  -5.96%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = "'2.43"
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = "'2.607.42"
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = "'360.80"
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = "'119.72"
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').Value = "'0.227"
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('E49').Value = '  -2.08%  '
$ws.Range('D50').Value = "'22.44"
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('D51').Value = "'1.93"
$ws.Range('E51').Value = '  -3.21%  '
